$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4 (ALC)
$ws.Range("H4").Value = 6700
$ws.Range("I4").Value = 6700
$ws.Range("K4").Value = 6700
$ws.Range("M4").Value = -6586

# Row 10 (ALC)
$ws.Range("H10").Value = 9000
$ws.Range("I10").Value = 9000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 9000
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -8707
$ws.Range("N10").ClearContents()

# Row 12 (ALC)
$ws.Range("H12").Value = 367.14285
$ws.Range("I12").Value = 300
$ws.Range("J12").Value = 394
$ws.Range("K12").Value = 300
$ws.Range("L12").Value = 394
$ws.Range("M12").Value = -130
$ws.Range("N12").Value = -734

# Row 61 (ALC)
$ws.Range("H61").Value = 2971.6667
$ws.Range("I61").Value = 1949.5
$ws.Range("K61").Value = 5848.5
$ws.Range("M61").Value = -5676.5

# Row 123 (ALC)
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Row 137 (ALC)
$ws.Range("H137").Value = 3434.8
$ws.Range("J137").Value = 3269.8
$ws.Range("L137").Value = 9809.400000000001
$ws.Range("N137").Value = -14909.4

# Row 138 (ALC)
$ws.Range("H138").Value = 3632
$ws.Range("I138").Value = 3163.3333
$ws.Range("J138").Value = 3749.1667
$ws.Range("K138").Value = 9489.999899999999
$ws.Range("L138").Value = 11247.5001
$ws.Range("M138").Value = -4349.999899999999
$ws.Range("N138").Value = -21527.5001

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws.Range("H2").Value = 1354
$ws.Range("I2").Value = 590.3333
$ws.Range("K2").Value = 590.3333
$ws.Range("M2").Value = -477.3333

# Row 5 (ARM)
$ws.Range("H5").Value = 673.3333
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 20
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 20
$ws.Range("M5").Value = -888
$ws.Range("N5").Value = -244

# Row 45 (ARM)
$ws.Range("H45").Value = 1850
$ws.Range("I45").Value = 1850
$ws.Range("K45").Value = 1850
$ws.Range("M45").Value = -1473

# Row 61 (ARM)
$ws.Range("H61").Value = 2958.1667
$ws.Range("I61").Value = 2958.1667
$ws.Range("K61").Value = 2958.1667
$ws.Range("M61").Value = -2746.1667

# Row 116 (ARM)
$ws.Range("H116").Value = 1354
$ws.Range("I116").Value = 590.3333
$ws.Range("K116").Value = 590.3333
$ws.Range("M116").Value = 1703.6667

# Row 128 (ARM)
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

# Row 132 (ARM)
$ws.Range("H132").Value = 5531.3335
$ws.Range("I132").Value = 594
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 1782
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = 748
$ws.Range("N132").Value = -29060

# Row 135 (ARM)
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# Row 136 (ARM)
$ws.Range("H136").Value = 2958.1667
$ws.Range("I136").Value = 2958.1667
$ws.Range("K136").Value = 8874.500100000001
$ws.Range("M136").Value = -6324.500100000001

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws.Range("H3").Value = 1354
$ws.Range("I3").Value = 590.3333
$ws.Range("K3").Value = 590.3333
$ws.Range("M3").Value = -476.3333

# Row 4 (BSM)
$ws.Range("H4").Value = 673.3333
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 20
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 20
$ws.Range("M4").Value = -885
$ws.Range("N4").Value = -250

# Row 80 (BSM)
$ws.Range("H80").Value = 593.2857
$ws.Range("I80").Value = 336.5
$ws.Range("K80").Value = 336.5
$ws.Range("M80").Value = 661.5

# Row 83 (BSM)
$ws.Range("H83").Value = 593.2857
$ws.Range("I83").Value = 336.5
$ws.Range("K83").Value = 1682.5
$ws.Range("M83").Value = 3309.5

# Row 86 (BSM)
$ws.Range("H86").Value = 1066.6666
$ws.Range("I86").Value = 1100
$ws.Range("J86").Value = 1025
$ws.Range("K86").Value = 1100
$ws.Range("L86").Value = 1025
$ws.Range("M86").Value = 23
$ws.Range("N86").Value = -3271

# Row 89 (BSM)
$ws.Range("H89").Value = 1066.6666
$ws.Range("I89").Value = 1100
$ws.Range("J89").Value = 1025
$ws.Range("K89").Value = 5500
$ws.Range("L89").Value = 5125
$ws.Range("M89").Value = 116
$ws.Range("N89").Value = -16357

# Row 105 (BSM)
$ws.Range("H105").Value = 3586.6667
$ws.Range("I105").Value = 3505
$ws.Range("K105").Value = 3505
$ws.Range("M105").Value = -1758

# Row 134 (BSM)
$ws.Range("H134").Value = 5475
$ws.Range("I134").Value = 3660
$ws.Range("K134").Value = 10980
$ws.Range("M134").Value = -8445

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (CRP)
$ws.Range("H7").Value = 81.933334
$ws.Range("I7").Value = 60.77778
$ws.Range("J7").Value = 113.666664
$ws.Range("K7").Value = 60.77778
$ws.Range("L7").Value = 113.666664
$ws.Range("M7").Value = 52.22222
$ws.Range("N7").Value = -339.666664

# Row 31 (CRP)
$ws.Range("H31").Value = 3718.04
$ws.Range("I31").Value = 1425.8889
$ws.Range("J31").Value = 5007.375
$ws.Range("K31").Value = 1425.8889
$ws.Range("L31").Value = 5007.375
$ws.Range("M31").Value = -1130.8889
$ws.Range("N31").Value = -5597.375

# Row 34 (CRP)
$ws.Range("H34").Value = 3718.04
$ws.Range("I34").Value = 1425.8889
$ws.Range("J34").Value = 5007.375
$ws.Range("K34").Value = 1425.8889
$ws.Range("L34").Value = 5007.375
$ws.Range("M34").Value = -1223.8889
$ws.Range("N34").Value = -5411.375

# Row 112 (CRP)
$ws.Range("H112").Value = 99980
$ws.Range("J112").Value = 99980
$ws.Range("L112").Value = 99980
$ws.Range("N112").Value = -102934

$ws = $wb.Worksheets.Item("CUL")
# Row 2 (CUL)
$ws.Range("H2").Value = 39.95652
$ws.Range("I2").Value = 17
$ws.Range("J2").Value = 57.615383
$ws.Range("K2").Value = 102
$ws.Range("L2").Value = 345.692298
$ws.Range("M2").Value = 11
$ws.Range("N2").Value = -571.6922979999999

# Row 38 (CUL)
$ws.Range("H38").Value = 50.125
$ws.Range("J38").Value = 65.8
$ws.Range("L38").Value = 197.4
$ws.Range("N38").Value = -891.4

# Row 68 (CUL)
$ws.Range("H68").Value = 900.25
$ws.Range("I68").Value = 1251
$ws.Range("J68").Value = 783.3333
$ws.Range("K68").Value = 3753
$ws.Range("L68").Value = 2349.9999
$ws.Range("M68").Value = -2942
$ws.Range("N68").Value = -3971.9999

# Row 71 (CUL)
$ws.Range("H71").Value = 900.25
$ws.Range("I71").Value = 1251
$ws.Range("J71").Value = 783.3333
$ws.Range("K71").Value = 11259
$ws.Range("L71").Value = 7049.9997
$ws.Range("M71").Value = -7203
$ws.Range("N71").Value = -15161.9997

# Row 80 (CUL)
$ws.Range("J80").Value = 6333.3335
$ws.Range("L80").Value = 19000.0005
$ws.Range("N80").Value = -20872.0005

# Row 83 (CUL)
$ws.Range("J83").Value = 6333.3335
$ws.Range("L83").Value = 57000.0015
$ws.Range("N83").Value = -66360.0015

# Row 94 (CUL)
$ws.Range("H94").Value = 5900
$ws.Range("I94").Value = 5900
$ws.Range("K94").Value = 17700
$ws.Range("M94").Value = -17024

$ws = $wb.Worksheets.Item("GSM")
# Row 97 (GSM)
$ws.Range("H97").Value = 3999.5
$ws.Range("I97").Value = 3999.5
$ws.Range("K97").Value = 3999.5
$ws.Range("M97").Value = -3503.5

# Row 132 (GSM)
$ws.Range("H132").Value = 2328.5715
$ws.Range("I132").Value = 2260
$ws.Range("K132").Value = 6780
$ws.Range("M132").Value = -4250

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Range("H7").Value = 8677.799999999999
$ws.Range("I7").Value = 8347.25
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 8347.25
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = -8235.25
$ws.Range("N7").Value = -10224

# Row 46 (LTW)
$ws.Range("H46").Value = 4948.421
$ws.Range("J46").Value = 4942.294
$ws.Range("L46").Value = 4942.294
$ws.Range("N46").Value = -5318.294

# Row 55 (LTW)
$ws.Range("H55").Value = 3660.5386
$ws.Range("I55").Value = 3082.1667
$ws.Range("J55").Value = 4156.2856
$ws.Range("K55").Value = 3082.1667
$ws.Range("L55").Value = 4156.2856
$ws.Range("M55").Value = -2909.1667
$ws.Range("N55").Value = -4502.2856

# Row 68 (LTW)
$ws.Range("H68").Value = 2637.5
$ws.Range("I68").Value = 2642.8572
$ws.Range("J68").Value = 2600
$ws.Range("K68").Value = 2642.8572
$ws.Range("L68").Value = 2600
$ws.Range("M68").Value = -1893.8572
$ws.Range("N68").Value = -4098

# Row 71 (LTW)
$ws.Range("H71").Value = 2637.5
$ws.Range("I71").Value = 2642.8572
$ws.Range("J71").Value = 2600
$ws.Range("K71").Value = 13214.286
$ws.Range("L71").Value = 13000
$ws.Range("M71").Value = -9470.286
$ws.Range("N71").Value = -20488

# Row 126 (LTW)
$ws.Range("H126").Value = 8677.799999999999
$ws.Range("I126").Value = 8347.25
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 25041.75
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -22571.75
$ws.Range("N126").Value = -34940

# Row 132 (LTW)
$ws.Range("H132").Value = 11333
$ws.Range("I132").Value = 9500
$ws.Range("J132").Value = 14999
$ws.Range("K132").Value = 28500
$ws.Range("L132").Value = 44997
$ws.Range("M132").Value = -25970
$ws.Range("N132").Value = -50057

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (WVR)
$ws.Range("H81").Value = 12689
$ws.Range("J81").Value = 1500
$ws.Range("L81").Value = 3000
$ws.Range("N81").Value = -5122

# Row 84 (WVR)
$ws.Range("H84").Value = 12689
$ws.Range("J84").Value = 1500
$ws.Range("L84").Value = 15000
$ws.Range("N84").Value = -25608

# Row 132 (WVR)
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 6000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -11060

# Row 133 (WVR)
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# Row 136 (WVR)
$ws.Range("H136").Value = 4850.375
$ws.Range("I136").Value = 4686.143
$ws.Range("K136").Value = 14058.429
$ws.Range("M136").Value = -11508.429
